$d = $word.ActiveDocument
$d.Content.Find.Execute("jmbg}", $true, $false, $false, $false, $false, $true, 1, $false, "jmbgNum}", 2)
